$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 47.63324766666667
$ws.Range("H2").Value = 142.899743
$ws.Range("I2").Value = 0.1341008591511219
$ws.Range("J2").Value = 0.1341008591511219
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 6940.199655925097
$ws.Range("R2").Value = 62461.79690332587
$ws.Range("S2").Value = 0.03843247145628267
$ws.Range("T2").Value = 0.03843247145628267
$ws.Range("G3").Value = 47.63324766666667
$ws.Range("H3").Value = 142.899743
$ws.Range("I3").Value = 0.1341008591511219
$ws.Range("J3").Value = 0.1341008591511219
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 8040.478043181029
$ws.Range("R3").Value = 72364.30238862925
$ws.Range("S3").Value = 0.04452543993105513
$ws.Range("T3").Value = 0.04452543993105514
$ws.Range("G4").Value = 47.63324766666667
$ws.Range("H4").Value = 142.899743
$ws.Range("I4").Value = 0.1341008591511219
$ws.Range("J4").Value = 0.1341008591511219
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 6103.064857814973
$ws.Range("R4").Value = 54927.58372033475
$ws.Range("S4").Value = 0.03379670291524928
$ws.Range("T4").Value = 0.03379670291524929
$ws.Range("G5").Value = 47.63324766666667
$ws.Range("H5").Value = 142.899743
$ws.Range("I5").Value = 0.1341008591511219
$ws.Range("J5").Value = 0.1341008591511219
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 3132.413762834232
$ws.Range("R5").Value = 28191.72386550809
$ws.Range("S5").Value = 0.01734624484853479
$ws.Range("T5").Value = 0.01734624484853479
$ws.Range("I6").Value = 0.4098937442001861
$ws.Range("J6").Value = 0.4098937442001861
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 21213.46902974098
$ws.Range("R6").Value = 190921.2212676689
$ws.Range("S6").Value = 0.1174729955035541
$ws.Range("T6").Value = 0.1174729955035541
$ws.Range("I7").Value = 0.4098937442001861
$ws.Range("J7").Value = 0.4098937442001861
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.1360968110199313
$ws.Range("T7").Value = 0.1360968110199313
$ws.Range("I8").Value = 0.4098937442001861
$ws.Range("J8").Value = 0.4098937442001861
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 18654.67620045019
$ws.Range("R8").Value = 167892.0858040517
$ws.Range("S8").Value = 0.1033032688026368
$ws.Range("T8").Value = 0.1033032688026368
$ws.Range("I9").Value = 0.4098937442001861
$ws.Range("J9").Value = 0.4098937442001861
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 9574.56062369736
$ws.Range("R9").Value = 86171.04561327623
$ws.Range("S9").Value = 0.0530206688740639
$ws.Range("T9").Value = 0.0530206688740639
$ws.Range("G10").Value = 11.826626
$ws.Range("H10").Value = 35.479878
$ws.Range("I10").Value = 0.03329524618093251
$ws.Range("J10").Value = 0.03329524618093251
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 1723.148215094161
$ws.Range("R10").Value = 15508.33393584745
$ws.Range("S10").Value = 0.009542210292900186
$ws.Range("T10").Value = 0.009542210292900187
$ws.Range("G11").Value = 11.826626
$ws.Range("H11").Value = 35.479878
$ws.Range("I11").Value = 0.03329524618093251
$ws.Range("J11").Value = 0.03329524618093251
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 1996.330952349869
$ws.Range("R11").Value = 17966.97857114883
$ws.Range("S11").Value = 0.01105500362341564
$ws.Range("T11").Value = 0.01105500362341565
$ws.Range("G12").Value = 11.826626
$ws.Range("H12").Value = 35.479878
$ws.Range("I12").Value = 0.03329524618093251
$ws.Range("J12").Value = 0.03329524618093251
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 1515.300112060821
$ws.Range("R12").Value = 13637.70100854739
$ws.Range("S12").Value = 0.008391218004046996
$ws.Range("T12").Value = 0.008391218004046999
$ws.Range("G13").Value = 11.826626
$ws.Range("H13").Value = 35.479878
$ws.Range("I13").Value = 0.03329524618093251
$ws.Range("J13").Value = 0.03329524618093251
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 777.7316866894539
$ws.Range("R13").Value = 6999.585180205086
$ws.Range("S13").Value = 0.004306814260569682
$ws.Range("T13").Value = 0.004306814260569683
$ws.Range("G14").Value = 150.1486076666667
$ws.Range("H14").Value = 450.445823
$ws.Range("I14").Value = 0.4227101504677595
$ws.Range("J14").Value = 0.4227101504677596
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 21876.76394769651
$ws.Range("R14").Value = 196890.8755292686
$ws.Range("S14").Value = 0.1211460977578473
$ws.Range("T14").Value = 0.1211460977578473
$ws.Range("G15").Value = 150.1486076666667
$ws.Range("H15").Value = 450.445823
$ws.Range("I15").Value = 0.4227101504677595
$ws.Range("J15").Value = 0.4227101504677596
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 25345.04032994732
$ws.Range("R15").Value = 228105.3629695259
$ws.Range("S15").Value = 0.140352235862182
$ws.Range("T15").Value = 0.140352235862182
$ws.Range("G16").Value = 150.1486076666667
$ws.Range("H16").Value = 450.445823
$ws.Range("I16").Value = 0.4227101504677595
$ws.Range("J16").Value = 0.4227101504677596
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 19237.96372888399
$ws.Range("R16").Value = 173141.6735599559
$ws.Range("S16").Value = 0.1065333172736774
$ws.Range("T16").Value = 0.1065333172736774
$ws.Range("G17").Value = 150.1486076666667
$ws.Range("H17").Value = 450.445823
$ws.Range("I17").Value = 0.4227101504677595
$ws.Range("J17").Value = 0.4227101504677596
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 9873.934450507671
$ws.Range("R17").Value = 88865.41005456905
$ws.Range("S17").Value = 0.05467849957405285
$ws.Range("T17").Value = 0.05467849957405286
